# Apply updated crypto price/volume data per the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '''61.976.19'
$ws.Cells.Item(2, 5).Value = '  -1.20%  '
$ws.Cells.Item(3, 4).Value = '''2.455.28'
$ws.Cells.Item(3, 5).Value = '  +0.72%  '
$ws.Cells.Item(4, 5).Value = '  -0.14%  '
$ws.Cells.Item(5, 4).Value = '''581.52'
$ws.Cells.Item(5, 5).Value = '  +0.42%  '
$ws.Cells.Item(6, 4).Value = '''141.98'
$ws.Cells.Item(6, 5).Value = '  -1.41%  '
$ws.Cells.Item(7, 5).Value = '  +0.02%  '
$ws.Cells.Item(8, 4).Value = '''0.530'
$ws.Cells.Item(8, 5).Value = '  +0.25%  '
$ws.Cells.Item(9, 4).Value = '''2.448.84'
$ws.Cells.Item(9, 5).Value = '  +0.54%  '
$ws.Cells.Item(10, 4).Value = '''0.110'
$ws.Cells.Item(10, 5).Value = '  +2.54%  '
$ws.Cells.Item(11, 4).Value = '''0.162'
$ws.Cells.Item(11, 5).Value = '  +3.02%  '
$ws.Cells.Item(12, 4).Value = '''5.17'
$ws.Cells.Item(12, 5).Value = '  -0.93%  '
$ws.Cells.Item(13, 4).Value = '''0.341'
$ws.Cells.Item(13, 5).Value = '  -2.08%  '
$ws.Cells.Item(14, 4).Value = '''25.97'
$ws.Cells.Item(14, 5).Value = '  -2.25%  '
$ws.Cells.Item(15, 2).Value = 'ShibaInu'
$ws.Cells.Item(15, 3).Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Cells.Item(15, 4).Value = '''0.0000173'
$ws.Cells.Item(15, 5).Value = '  -0.57%  '
$ws.Cells.Item(16, 2).Value = 'WrappedliquidstakedEther2.0'
$ws.Cells.Item(16, 3).Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Cells.Item(16, 4).Value = '''2.893.75'
$ws.Cells.Item(16, 5).Value = '  +0.41%  '
$ws.Cells.Item(17, 4).Value = '''61.871.37'
$ws.Cells.Item(17, 5).Value = '  -1.17%  '
$ws.Cells.Item(18, 4).Value = '''2.441.53'
$ws.Cells.Item(18, 5).Value = '  +0.06%  '
$ws.Cells.Item(19, 4).Value = '''10.66'
$ws.Cells.Item(19, 5).Value = '  -3.48%  '
$ws.Cells.Item(20, 4).Value = '''7.24'
$ws.Cells.Item(20, 5).Value = '  +1.54%  '
$ws.Cells.Item(21, 4).Value = '''325.87'
$ws.Cells.Item(21, 5).Value = '  -1.78%  '
$ws.Cells.Item(22, 2).Value = 'Polkadot'
$ws.Cells.Item(22, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Cells.Item(22, 4).Value = '''4.09'
$ws.Cells.Item(22, 5).Value = '  -1.16%  '
$ws.Cells.Item(23, 2).Value = 'SuiNetwork'
$ws.Cells.Item(23, 3).Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Cells.Item(23, 4).Value = '''1.94'
$ws.Cells.Item(23, 5).Value = '  -2.74%  '
$ws.Cells.Item(25, 2).Value = 'Litecoin'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Cells.Item(25, 4).Value = '''65.14'
$ws.Cells.Item(25, 5).Value = '  -1.28%  '
$ws.Cells.Item(26, 2).Value = 'Aptos'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Cells.Item(26, 4).Value = '''9.17'
$ws.Cells.Item(26, 5).Value = '  +1.79%  '
$ws.Cells.Item(27, 2).Value = 'Bittensor'
$ws.Cells.Item(27, 3).Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Cells.Item(27, 4).Value = '''586.21'
$ws.Cells.Item(27, 5).Value = '  -8.31%  '
$ws.Cells.Item(28, 2).Value = 'WrappedeETH'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Cells.Item(28, 4).Value = '''2.569.32'
$ws.Cells.Item(28, 5).Value = '  +0.26%  '
$ws.Cells.Item(29, 2).Value = 'Binance-PegBSC-USD'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Cells.Item(29, 4).Value = '''0.999'
$ws.Cells.Item(29, 5).Value = '  -0.04%  '
$ws.Cells.Item(30, 2).Value = 'PEPE'
$ws.Cells.Item(30, 3).Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Cells.Item(30, 4).Value = '0.0₃0938'
$ws.Cells.Item(30, 5).Value = '  -2.57%  '
$ws.Cells.Item(31, 2).Value = 'InternetComputer(DFINITY)'
$ws.Cells.Item(31, 3).Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Cells.Item(31, 4).Value = '''7.94'
$ws.Cells.Item(31, 5).Value = '  -1.74%  '
$ws.Cells.Item(32, 2).Value = 'Fetch.AI'
$ws.Cells.Item(32, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Cells.Item(32, 4).Value = '''1.38'
$ws.Cells.Item(32, 5).Value = '  -4.41%  '
$ws.Cells.Item(33, 2).Value = 'PancakeSwap'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Cells.Item(33, 4).Value = '''1.88'
$ws.Cells.Item(33, 5).Value = '  -0.68%  '
$ws.Cells.Item(34, 2).Value = 'Kaspa'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Cells.Item(34, 4).Value = '''0.133'
$ws.Cells.Item(34, 5).Value = '  -4.02%  '
$ws.Cells.Item(35, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Cells.Item(35, 4).Value = '''1.00'
$ws.Cells.Item(35, 5).Value = '  -0.08%  '
$ws.Cells.Item(36, 2).Value = 'NEARProtocol'
$ws.Cells.Item(36, 3).Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Cells.Item(36, 4).Value = '''4.77'
$ws.Cells.Item(36, 5).Value = '  -4.90%  '
$ws.Cells.Item(37, 2).Value = 'Monero'
$ws.Cells.Item(37, 3).Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Cells.Item(37, 4).Value = '''152.63'
$ws.Cells.Item(37, 5).Value = '  +2.74%  '
$ws.Cells.Item(38, 4).Value = '''0.372'
$ws.Cells.Item(38, 5).Value = '  -1.08%  '
$ws.Cells.Item(39, 2).Value = 'ImmutableX'
$ws.Cells.Item(39, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Cells.Item(39, 4).Value = '''1.41'
$ws.Cells.Item(39, 5).Value = '  -3.05%  '
$ws.Cells.Item(40, 2).Value = 'EthereumClassic'
$ws.Cells.Item(40, 3).Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Cells.Item(40, 4).Value = '''18.31'
$ws.Cells.Item(40, 5).Value = '  -1.03%  '
$ws.Cells.Item(41, 2).Value = 'RenderToken'
$ws.Cells.Item(41, 3).Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
$ws.Cells.Item(41, 4).Value = '''5.18'
$ws.Cells.Item(41, 5).Value = '  -2.29%  '
$ws.Cells.Item(42, 2).Value = 'USDe'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Cells.Item(42, 4).Value = '''0.999'
$ws.Cells.Item(42, 5).Value = '  +0.00%  '
$ws.Cells.Item(43, 2).Value = 'OKB'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Cells.Item(43, 4).Value = '''42.11'
$ws.Cells.Item(43, 5).Value = '  -1.04%  '
$ws.Cells.Item(44, 2).Value = 'Stacks'
$ws.Cells.Item(44, 3).Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Cells.Item(44, 4).Value = '''1.68'
$ws.Cells.Item(44, 5).Value = '  -4.15%  '
$ws.Cells.Item(45, 2).Value = 'dogwifhat'
$ws.Cells.Item(45, 3).Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Cells.Item(45, 4).Value = '''2.37'
$ws.Cells.Item(45, 5).Value = '  -4.82%  '
$ws.Cells.Item(46, 2).Value = 'BabyDogeCoin'
$ws.Cells.Item(46, 3).Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Cells.Item(46, 4).Value = '0.0₆0282'
$ws.Cells.Item(46, 5).Value = '  +16.50%  '
$ws.Cells.Item(47, 2).Value = 'Aave'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Cells.Item(47, 4).Value = '''141.60'
$ws.Cells.Item(47, 5).Value = '  -1.86%  '
$ws.Cells.Item(48, 2).Value = 'Filecoin'
$ws.Cells.Item(48, 3).Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Cells.Item(48, 4).Value = '''3.57'
$ws.Cells.Item(48, 5).Value = '  -3.29%  '
$ws.Cells.Item(49, 2).Value = 'Mantle'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Cells.Item(49, 4).Value = '''0.600'
$ws.Cells.Item(49, 5).Value = '  +0.23%  '
$ws.Cells.Item(50, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(50, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Cells.Item(50, 4).Value = '''19.75'
$ws.Cells.Item(50, 5).Value = '  +0.44%  '
$ws.Cells.Item(51, 4).Value = '''0.0512'
$ws.Cells.Item(51, 5).Value = '  -2.05%  '

Write-Host "Updated cryptos list"
